# Updated cryptos list on Mon Sep 16 07:35:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text
# formatting instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.896.87"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").Value = "2.298.82"
$ws.Range("E3").Value = "  -5.22%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "548.21"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").Value = "131.12"
$ws.Range("E6").Value = "  -4.50%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -2.90%  "

$ws.Range("D9").Value = "2.296.47"
$ws.Range("E9").Value = "  -5.25%  "

$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  -2.69%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("E13").Value = "  -5.12%  "

$ws.Range("D14").Value = "23.98"
$ws.Range("E14").Value = "  -3.82%  "

$ws.Range("D15").Value = "2.704.58"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").Value = "58.839.89"
$ws.Range("E16").Value = "  -2.13%  "

$ws.Range("E17").Value = "  -3.42%  "

$ws.Range("D18").Value = "2.268.02"
$ws.Range("E18").Value = "  -5.83%  "

$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  -5.20%  "

$ws.Range("D20").Value = "4.32"
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").Value = "315.77"

$ws.Range("D22").Value = "6.49"
$ws.Range("E22").Value = "  -4.01%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "63.17"
$ws.Range("E24").Value = "  -3.29%  "

$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "8.11"
$ws.Range("E27").Value = "  -6.42%  "

$ws.Range("E28").Value = "  -7.82%  "

$ws.Range("D29").Value = "1.76"
$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").Value = "169.72"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Value = "0.0₃0732"
$ws.Range("E31").Value = "  -5.73%  "

$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  +3.65%  "

$ws.Range("D33").Value = "5.82"
$ws.Range("E33").Value = "  -4.91%  "

$ws.Range("E34").Value = "  -4.77%  "

$ws.Range("D36").Value = "17.79"
$ws.Range("E36").Value = "  -4.14%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -6.46%  "

$ws.Range("D39").Value = "3.98"
$ws.Range("E39").Value = "  -5.95%  "

$ws.Range("D40").Value = "38.02"
$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("E41").Value = "  -5.38%  "

$ws.Range("D42").Value = "303.51"
$ws.Range("E42").Value = "  -7.10%  "

$ws.Range("D43").Value = "140.13"

$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").Value = "0.0501"
$ws.Range("E46").Value = "  -3.07%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "18.74"
$ws.Range("E47").Value = "  -5.80%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.558"
$ws.Range("E48").Value = "  -3.25%  "

$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").Value = "16.72"
$ws.Range("E50").Value = "  -4.56%  "

$ws.Range("E51").Value = "  -0.32%  "
